# Horarios actualizados Línea 141 - 618
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with the
# latest scrape: refreshed "Última actualización" / "Total filas" headers,
# a handful of rows that got re-sorted because newer scrape timestamps (col A)
# landed on an already-seen arrival time (col B), a new row inserted in the
# middle of the LP1912 sheet, and new rows appended at the tail of LP1912 and
# LP1912-215.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 22:04:33"
$ws1.Range("A3").Value = "Total filas: 377"

# --- Re-sorted rows: newer scrape (col A) reshuffled rows sharing the same
#     arrival time (col B). Only columns A, C, D change; B and E stay put.
$ws1.Cells.Item(62,1).Value = "07:49:32"
$ws1.Cells.Item(62,3).Value = "14_ABASTO"
$ws1.Cells.Item(62,4).Value = 88

$ws1.Cells.Item(63,1).Value = "08:38:24"
$ws1.Cells.Item(63,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(63,4).Value = 39

$ws1.Cells.Item(86,1).Value = "10:04:30"
$ws1.Cells.Item(86,3).Value = "14_ABASTO"
$ws1.Cells.Item(86,4).Value = 25

$ws1.Cells.Item(87,1).Value = "08:38:24"
$ws1.Cells.Item(87,3).Value = "15_ABASTO"
$ws1.Cells.Item(87,4).Value = 111

$ws1.Cells.Item(133,1).Value = "11:33:52"
$ws1.Cells.Item(133,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(133,4).Value = 59

$ws1.Cells.Item(134,1).Value = "10:36:50"
$ws1.Cells.Item(134,3).Value = "14_ABASTO"
$ws1.Cells.Item(134,4).Value = 116

$ws1.Cells.Item(158,1).Value = "12:11:21"
$ws1.Cells.Item(158,3).Value = "14_ABASTO"
$ws1.Cells.Item(158,4).Value = 81

$ws1.Cells.Item(159,1).Value = "11:53:44"
$ws1.Cells.Item(159,3).Value = "215A_EL PATO"
$ws1.Cells.Item(159,4).Value = 99

$ws1.Cells.Item(169,1).Value = "13:14:31"
$ws1.Cells.Item(169,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(169,4).Value = 48

$ws1.Cells.Item(171,1).Value = "12:33:02"
$ws1.Cells.Item(171,3).Value = "10_OLMOS"
$ws1.Cells.Item(171,4).Value = 89

$ws1.Cells.Item(204,1).Value = "14:32:44"
$ws1.Cells.Item(204,3).Value = "10_OLMOS"
$ws1.Cells.Item(204,4).Value = 81

$ws1.Cells.Item(206,1).Value = "13:55:43"
$ws1.Cells.Item(206,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(206,4).Value = 118

$ws1.Cells.Item(276,1).Value = "17:13:30"
$ws1.Cells.Item(276,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(276,4).Value = 76

$ws1.Cells.Item(277,1).Value = "16:37:37"
$ws1.Cells.Item(277,3).Value = "14_ABASTO"
$ws1.Cells.Item(277,4).Value = 112

$ws1.Cells.Item(358,1).Value = "20:46:15"
$ws1.Cells.Item(358,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(358,4).Value = 48

$ws1.Cells.Item(359,1).Value = "19:35:34"
$ws1.Cells.Item(359,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(359,4).Value = 119

# --- A brand-new scrape row lands right before the old row 370, so every row
#     from 370-376 shifts down by one to 371-377.
$ws1.Rows.Item(370).Insert()

$ws1.Cells.Item(370,1).Value = "22:04:33"
$ws1.Cells.Item(370,2).Value = "22:08"
$ws1.Cells.Item(370,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(370,4).Value = 4
$ws1.Cells.Item(370,5).Value = "LP1912"

# --- Six more brand-new rows appended at the tail (377 was the old last row,
#     now at 377 after the insert above; new data continues at 378-382).
$ws1.Cells.Item(378,1).Value = "22:04:33"
$ws1.Cells.Item(378,2).Value = "23:04"
$ws1.Cells.Item(378,3).Value = "15_ABASTO"
$ws1.Cells.Item(378,4).Value = 60
$ws1.Cells.Item(378,5).Value = "LP1912"

$ws1.Cells.Item(379,1).Value = "22:04:33"
$ws1.Cells.Item(379,2).Value = "23:22"
$ws1.Cells.Item(379,3).Value = "14_ABASTO"
$ws1.Cells.Item(379,4).Value = 78
$ws1.Cells.Item(379,5).Value = "LP1912"

$ws1.Cells.Item(380,1).Value = "22:04:33"
$ws1.Cells.Item(380,2).Value = "23:34"
$ws1.Cells.Item(380,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(380,4).Value = 90
$ws1.Cells.Item(380,5).Value = "LP1912"

$ws1.Cells.Item(381,1).Value = "22:04:33"
$ws1.Cells.Item(381,2).Value = "23:40"
$ws1.Cells.Item(381,3).Value = "215A_EL PATO"
$ws1.Cells.Item(381,4).Value = 96
$ws1.Cells.Item(381,5).Value = "LP1912"

$ws1.Cells.Item(382,1).Value = "22:04:33"
$ws1.Cells.Item(382,2).Value = "23:58"
$ws1.Cells.Item(382,3).Value = "11X44_ETCHEVERRY"
$ws1.Cells.Item(382,4).Value = 114
$ws1.Cells.Item(382,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 22:04:33"
$ws2.Range("A3").Value = "Total filas: 57"

$ws2.Cells.Item(62,1).Value = "22:04:33"
$ws2.Cells.Item(62,2).Value = "23:40"
$ws2.Cells.Item(62,3).Value = "215A_EL PATO"
$ws2.Cells.Item(62,4).Value = 96
$ws2.Cells.Item(62,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 22:04:33"
